$wb = $excel.ActiveWorkbook

# Scheduled market-data refresh for the per-job "Profits" sheets (ALC, ARM,
# BSM, CRP, CUL, GSM, LTW, WVR). Columns H-N hold price/profit figures
# (currentAveragePrice*, LevePrice*, LeveProfit*) refreshed by the runner.
# A few rows gain or lose their LeveProfit* cells entirely depending on
# whether a profit figure is computable from the refreshed prices.

# ----- ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 1812.6207
$ws.Range("I86").Value = 791.8946999999999
$ws.Range("J86").Value = 3752
$ws.Range("K86").Value = 791.8946999999999
$ws.Range("L86").Value = 3752
$ws.Range("M86").Value = 331.1053000000001
$ws.Range("N86").Value = -5998
$ws.Range("H88").Value = 1224
$ws.Range("I88").Value = 369.2
$ws.Range("J88").Value = 1758.25
$ws.Range("K88").Value = 369.2
$ws.Range("L88").Value = 1758.25
$ws.Range("M88").Value = 36.80000000000001
$ws.Range("N88").Value = -2570.25
$ws.Range("H89").Value = 1812.6207
$ws.Range("I89").Value = 791.8946999999999
$ws.Range("J89").Value = 3752
$ws.Range("K89").Value = 3959.4735
$ws.Range("L89").Value = 18760
$ws.Range("M89").Value = 1656.5265
$ws.Range("N89").Value = -29992
$ws.Range("H91").Value = 1224
$ws.Range("I91").Value = 369.2
$ws.Range("J91").Value = 1758.25
$ws.Range("K91").Value = 369.2
$ws.Range("L91").Value = 1758.25
$ws.Range("M91").Value = 1034.8
$ws.Range("N91").Value = -4566.25
$ws.Range("H113").Value = 1983.3334
$ws.Range("I113").Value = 2320
$ws.Range("J113").Value = 1742.8572
$ws.Range("K113").Value = 2320
$ws.Range("L113").Value = 1742.8572
$ws.Range("M113").Value = 934
$ws.Range("N113").Value = -8250.8572
$ws.Range("H132").Value = 3660.5356
$ws.Range("I132").Value = 3781.2964
$ws.Range("K132").Value = 11343.8892
$ws.Range("M132").Value = -8813.889200000001
$ws.Range("H137").Value = 27029964
$ws.Range("I137").Value = 1589.25
$ws.Range("J137").Value = 58828052
$ws.Range("K137").Value = 4767.75
$ws.Range("L137").Value = 176484156
$ws.Range("M137").Value = -2217.75
$ws.Range("N137").Value = -176489256
$ws.Range("H138").Value = 3246.423
$ws.Range("I138").Value = 2715.2104
$ws.Range("J138").Value = 3751.075
$ws.Range("K138").Value = 8145.6312
$ws.Range("L138").Value = 11253.225
$ws.Range("M138").Value = -3005.6312
$ws.Range("N138").Value = -21533.225

# ----- ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H8").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("N8").ClearContents()
$ws.Range("H122").Value = 1588.8462
$ws.Range("I122").Value = 1457.4286
$ws.Range("J122").Value = 1637.2632
$ws.Range("K122").Value = 4372.2858
$ws.Range("L122").Value = 4911.7896
$ws.Range("M122").Value = -1922.2858
$ws.Range("N122").Value = -9811.7896
$ws.Range("H128").Value = 46249.8
$ws.Range("J128").Value = 46249.8
$ws.Range("L128").Value = 46249.8
$ws.Range("N128").Value = -56209.8

# ----- BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1769.1852
$ws.Range("I86").Value = 1634.7
$ws.Range("K86").Value = 1634.7
$ws.Range("M86").Value = -511.7
$ws.Range("H89").Value = 1769.1852
$ws.Range("I89").Value = 1634.7
$ws.Range("K89").Value = 8173.5
$ws.Range("M89").Value = -2557.5
$ws.Range("H107").Value = 492459.25
$ws.Range("I107").Value = 833936.4
$ws.Range("J107").Value = 8700
$ws.Range("K107").Value = 833936.4
$ws.Range("L107").Value = 8700
$ws.Range("M107").Value = -832016.4
$ws.Range("N107").Value = -12540

# ----- CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1450
$ws.Range("J16").Value = 1450
$ws.Range("L16").Value = 1450
$ws.Range("N16").Value = -2024
$ws.Range("H113").Value = 1450
$ws.Range("J113").Value = 1450
$ws.Range("L113").Value = 1450
$ws.Range("N113").Value = -5790
$ws.Range("H132").Value = 3335.5715
$ws.Range("I132").Value = 3264.5454
$ws.Range("K132").Value = 9793.636200000001
$ws.Range("M132").Value = -7263.636200000001
$ws.Range("H134").Value = 3112.5
$ws.Range("I134").Value = 3019.125
$ws.Range("J134").Value = 3261.9
$ws.Range("K134").Value = 9057.375
$ws.Range("L134").Value = 9785.700000000001
$ws.Range("M134").Value = -6522.375
$ws.Range("N134").Value = -14855.7

# ----- CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 12072911
$ws.Range("I80").Value = 18108692
$ws.Range("J80").Value = 1350
$ws.Range("K80").Value = 54326076
$ws.Range("L80").Value = 4050
$ws.Range("M80").Value = -54325140
$ws.Range("N80").Value = -5922
$ws.Range("H83").Value = 12072911
$ws.Range("I83").Value = 18108692
$ws.Range("J83").Value = 1350
$ws.Range("K83").Value = 162978228
$ws.Range("L83").Value = 12150
$ws.Range("M83").Value = -162973548
$ws.Range("N83").Value = -21510
$ws.Range("H94").Value = 2000
$ws.Range("I94").Value = 1000
$ws.Range("J94").Value = 3000
$ws.Range("K94").Value = 3000
$ws.Range("L94").Value = 9000
$ws.Range("M94").Value = -2324
$ws.Range("N94").Value = -10352
$ws.Range("H103").Value = 8501256
$ws.Range("I103").Value = 17000512
$ws.Range("J103").Value = 2000
$ws.Range("K103").Value = 51001536
$ws.Range("L103").Value = 6000
$ws.Range("M103").Value = -51000657
$ws.Range("N103").Value = -7758
$ws.Range("H108").Value = 300
$ws.Range("I108").Value = 300
$ws.Range("K108").Value = 900
$ws.Range("M108").Value = 1980
$ws.Range("H114").Value = 1398.24
$ws.Range("I114").Value = 937.7
$ws.Range("J114").Value = 1705.2667
$ws.Range("K114").Value = 2813.1
$ws.Range("L114").Value = 5115.800099999999
$ws.Range("M114").Value = 440.8999999999996
$ws.Range("N114").Value = -11623.8001
$ws.Range("H122").Value = 7003906.5
$ws.Range("I122").Value = 13333911
$ws.Range("J122").Value = 917364.4
$ws.Range("K122").Value = 120005199
$ws.Range("L122").Value = 8256279.600000001
$ws.Range("M122").Value = -120002749
$ws.Range("N122").Value = -8261179.600000001

# ----- GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H6").Value = 166.66667
$ws.Range("I6").Value = 100
$ws.Range("J6").Value = 200
$ws.Range("K6").Value = 100
$ws.Range("L6").Value = 200
$ws.Range("M6").Value = 13
$ws.Range("N6").Value = -426
$ws.Range("H16").Value = 166.66667
$ws.Range("I16").Value = 100
$ws.Range("J16").Value = 200
$ws.Range("K16").Value = 100
$ws.Range("L16").Value = 200
$ws.Range("M16").Value = 150
$ws.Range("N16").Value = -700
$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()
$ws.Range("H102").Value = 1230.1666
$ws.Range("I102").Value = 1151.7142
$ws.Range("K102").Value = 1151.7142
$ws.Range("M102").Value = 470.2858000000001
$ws.Range("H107").Value = 416.6842
$ws.Range("I107").Value = 339.13333
$ws.Range("J107").Value = 707.5
$ws.Range("K107").Value = 339.13333
$ws.Range("L107").Value = 707.5
$ws.Range("M107").Value = 1580.86667
$ws.Range("N107").Value = -4547.5

# ----- LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H94").Value = 25000
$ws.Range("J94").Value = 25000
$ws.Range("L94").Value = 25000
$ws.Range("N94").Value = -26352
$ws.Range("H98").Value = 50000
$ws.Range("J98").Value = 50000
$ws.Range("L98").Value = 50000
$ws.Range("N98").Value = -55990
$ws.Range("H140").Value = 147751.4
$ws.Range("J140").Value = 147751.4
$ws.Range("L140").Value = 147751.4
$ws.Range("N140").Value = -158111.4

# ----- WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3503.9285
$ws.Range("I62").Value = 4083.6667
$ws.Range("J62").Value = 3069.125
$ws.Range("K62").Value = 4083.6667
$ws.Range("L62").Value = 3069.125
$ws.Range("M62").Value = -3459.6667
$ws.Range("N62").Value = -4317.125
$ws.Range("H65").Value = 3503.9285
$ws.Range("I65").Value = 4083.6667
$ws.Range("J65").Value = 3069.125
$ws.Range("K65").Value = 20418.3335
$ws.Range("L65").Value = 15345.625
$ws.Range("M65").Value = -17298.3335
$ws.Range("N65").Value = -21585.625
$ws.Range("H107").Value = 587.86365
$ws.Range("I107").Value = 413.70587
$ws.Range("J107").Value = 1180
$ws.Range("K107").Value = 1241.11761
$ws.Range("L107").Value = 3540
$ws.Range("M107").Value = 678.88239
$ws.Range("N107").Value = -7380
